$d = $word.ActiveDocument

function Replace-ParagraphXml {
    param($Paragraph, $NewParaXml)
    $rng = $Paragraph.Range
    $pkg = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $NewParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

# NOTE: this runtime's -like/-match/-eq operators are NOT case sensitive
# (even the "c"-prefixed variants), so we use the case-sensitive .NET
# [string]::Contains() method to find anchors unambiguously.
function Find-ParagraphIndexContaining {
    param($Anchor)
    $idx = 0
    $foundIdx = -1
    $matchCount = 0
    foreach ($p in $d.Paragraphs) {
        $idx = $idx + 1
        if ($p.Range.Text.Contains($Anchor)) {
            $matchCount = $matchCount + 1
            if ($foundIdx -lt 0) { $foundIdx = $idx }
        }
    }
    if ($matchCount -ne 1) {
        throw ("Anchor '" + $Anchor + "' matched " + $matchCount + " paragraphs (expected exactly 1)")
    }
    return $foundIdx
}

# ------------------------------------------------------------------
# Edit 4 first (remove the old "_GoBack" bookmark paragraph), BEFORE
# edit 2 inserts a new "_GoBack" bookmark elsewhere, so our anchor
# text lookup below is unambiguous.
# ------------------------------------------------------------------
$idxCela = Find-ParagraphIndexContaining "Cela permet une gestion simple"
$p4 = $d.Paragraphs($idxCela + 1)
if ($p4.Range.Text -ne [char]13) { throw "Paragraph after 'Cela permet...' was not the expected empty bookmark paragraph (text='" + $p4.Range.Text + "')" }
$target4 = '<w:p w14:paraId="50319C7D" w14:textId="1A4D5795" w:rsidR="00FB220F" w:rsidRDefault="00FB220F" w:rsidP="00FB220F"></w:p>'
Replace-ParagraphXml $p4 $target4

# --- Edit 1: "La table Employe garde les informations de chaque employés..." paragraph ---
$idx1 = Find-ParagraphIndexContaining "Stocker les employ"
$p1 = $d.Paragraphs($idx1)
$target1 = '<w:p w14:paraId="4A98154F" w14:textId="787B7574" w:rsidR="004A3882" w:rsidRDefault="004A3882" w:rsidP="004A3882"><w:r><w:tab/></w:r><w:r w:rsidRPr="004A3882"><w:rPr><w:b/></w:rPr><w:t>Stocker les employés</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> : </w:t></w:r><w:r><w:t xml:space="preserve">La table </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>Employe</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">garde les informations de chaque </w:t></w:r><w:r><w:t>employé</w:t></w:r><w:r><w:t xml:space="preserve">, comme leur prénom, nom email ou encore mot de passe. Leur </w:t></w:r><w:r><w:tab/><w:t xml:space="preserve">statut admin y est aussi défini, et chaque employé est lié à une ligue grâce à la clé étrangère ID</w:t></w:r><w:r><w:t>_</w:t></w:r><w:r><w:t>ligue</w:t></w:r><w:r w:rsidR="00677C99"><w:t>.</w:t></w:r></w:p>'
Replace-ParagraphXml $p1 $target1

# --- Edit 2: "Modification  du Code en Java :" heading paragraph ---
$idx2 = Find-ParagraphIndexContaining "Modification"
$p2 = $d.Paragraphs($idx2)
$target2 = '<w:p w14:paraId="70C90A89" w14:textId="59A95D61" w:rsidR="00FB220F" w:rsidRDefault="00FB220F" w:rsidP="00FB220F"><w:pPr><w:pStyle w:val="Titre2"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>Modification du</w:t></w:r><w:r><w:t xml:space="preserve"> Code en Java :</w:t></w:r></w:p>'
Replace-ParagraphXml $p2 $target2

# --- Edit 3: "Le code LigueConsole propose..." paragraph ---
$idx3 = Find-ParagraphIndexContaining "Le code "
$p3 = $d.Paragraphs($idx3)
$target3 = '<w:p w14:paraId="677E912C" w14:textId="77777777" w:rsidR="00FB220F" w:rsidRDefault="00FB220F" w:rsidP="00FB220F"><w:r><w:t>Le code LigueConsole propose une interface en ligne qui permet d’interagir avec l’utilisateur. Celui-ci permet de gérer des ligues et leurs employés. Voici les principales fonctionnalités :</w:t></w:r></w:p>'
Replace-ParagraphXml $p3 $target3

"DONE"
